$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

# Data to sync: rows 61-69 move from "Pending" to "Completed" with results filled in.
# Columns: L=Status, M=Result, N=Resultado_Real, O=Profit, P=ROI, Q=Enviado
$timestamp = "2025-09-01 04:32:43"

$updates = @(
    @{ Row = 61; Result = "Away Win"; Resultado = "Fallo";   Profit = -3.5; Roi = -100 },
    @{ Row = 62; Result = "Draw";     Resultado = "Fallo";   Profit = -5.3; Roi = -100 },
    @{ Row = 63; Result = "Home Win"; Resultado = "Fallo";   Profit = -1.4; Roi = -100 },
    @{ Row = 64; Result = "Draw";     Resultado = "Fallo";   Profit = -4;   Roi = -100 },
    @{ Row = 65; Result = "Home Win"; Resultado = "Acierto"; Profit = 2.65; Roi = 50   },
    @{ Row = 66; Result = "Home Win"; Resultado = "Fallo";   Profit = -0.4; Roi = -100 },
    @{ Row = 67; Result = "Draw";     Resultado = "Fallo";   Profit = -5.3; Roi = -100 },
    @{ Row = 68; Result = "Draw";     Resultado = "Fallo";   Profit = -4.8; Roi = -100 },
    @{ Row = 69; Result = "Draw";     Resultado = "Fallo";   Profit = -2.8; Roi = -100 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 12).Value = "Completed"      # L - Status
    $ws.Cells.Item($r, 13).Value = $u.Result         # M - Result
    $ws.Cells.Item($r, 14).Value = $u.Resultado      # N - Resultado_Real
    $ws.Cells.Item($r, 15).Value = $u.Profit         # O - Profit
    $ws.Cells.Item($r, 16).Value = $u.Roi            # P - ROI
    $ws.Cells.Item($r, 17).Value = $timestamp        # Q - Enviado
}
